$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> updated Coin / Link / Price / Volume(1h) values, per latest cryptos.xlsx refresh
$rows = @(
    @{ Row = 2; Coin = 'Bitcoin'; Link = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'; Price = '28.448.52'; Volume = '  +0.32%  ' }
    @{ Row = 3; Coin = 'Ethereum'; Link = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'; Price = '1.867.71'; Volume = '  +0.13%  ' }
    @{ Row = 4; Coin = 'TetherUSD'; Link = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'; Price = '1.007'; Volume = '  -1.19%  ' }
    @{ Row = 5; Coin = 'BNB'; Link = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; Price = '314.21'; Volume = '  -0.77%  ' }
    @{ Row = 6; Coin = 'USDC'; Link = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; Price = '1.005'; Volume = '  -1.24%  ' }
    @{ Row = 7; Coin = 'XRP'; Link = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; Price = '0.5066'; Volume = '  -0.63%  ' }
    @{ Row = 8; Coin = 'Cardano'; Link = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; Price = '0.3905'; Volume = '  -1.26%  ' }
    @{ Row = 9; Coin = 'Dogecoin'; Link = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; Price = '0.08319'; Volume = '  -0.17%  ' }
    @{ Row = 10; Coin = 'Polygon'; Link = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; Price = '1.103'; Volume = '  -0.27%  ' }
    @{ Row = 11; Coin = 'Polkadot'; Link = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'; Price = '6.170'; Volume = '  -0.86%  ' }
    @{ Row = 12; Coin = 'WrappedEther'; Link = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; Price = '1.861.19'; Volume = '  +0.59%  ' }
    @{ Row = 13; Coin = 'Solana'; Link = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; Price = '20.29'; Volume = '  -0.63%  ' }
    @{ Row = 14; Coin = 'Chainlink'; Link = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; Price = '7.226'; Volume = '  +0.55%  ' }
    @{ Row = 15; Coin = 'BinanceUSD'; Link = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; Price = '1.008'; Volume = '  -1.20%  ' }
    @{ Row = 16; Coin = 'ShibaInu'; Link = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; Price = '0.00001097'; Volume = '  -0.67%  ' }
    @{ Row = 17; Coin = 'Litecoin'; Link = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; Price = '91.13'; Volume = '  +0.56%  ' }
    @{ Row = 18; Coin = 'TRON'; Link = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; Price = '0.06715'; Volume = '  -0.37%  ' }
    @{ Row = 19; Coin = 'Avalanche'; Link = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; Price = '17.59'; Volume = '  -0.31%  ' }
    @{ Row = 20; Coin = 'Dai'; Link = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'; Price = '1.006'; Volume = '  -1.19%  ' }
    @{ Row = 21; Coin = 'Uniswap'; Link = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; Price = '5.899'; Volume = '  -0.65%  ' }
    @{ Row = 22; Coin = 'WrappedBTC'; Link = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; Price = '28.464.77'; Volume = '  +0.24%  ' }
    @{ Row = 23; Coin = 'Cosmos'; Link = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; Price = '11.03'; Volume = '  -0.75%  ' }
    @{ Row = 24; Coin = 'Toncoin'; Link = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; Price = '2.193'; Volume = '  -3.84%  ' }
    @{ Row = 25; Coin = 'WrappedliquidstakedEther2.0'; Link = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; Price = '2.072.77'; Volume = '  +0.86%  ' }
    @{ Row = 26; Coin = 'Monero'; Link = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; Price = '157.77'; Volume = '  -2.34%  ' }
    @{ Row = 27; Coin = 'EthereumClassic'; Link = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; Price = '20.50'; Volume = '  -0.63%  ' }
    @{ Row = 28; Coin = 'LidoDAOToken'; Link = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; Price = '2.409'; Volume = '  +2.23%  ' }
    @{ Row = 29; Coin = 'BitcoinCash'; Link = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; Price = '126.30'; Volume = '  -0.63%  ' }
    @{ Row = 30; Coin = 'Stellar'; Link = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; Price = '0.1033'; Volume = '  -1.27%  ' }
    @{ Row = 31; Coin = 'ImmutableX'; Link = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; Price = '1.034'; Volume = '  +0.33%  ' }
    @{ Row = 32; Coin = 'Filecoin'; Link = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; Price = '5.769'; Volume = '  +0.07%  ' }
    @{ Row = 33; Coin = 'HuobiToken'; Link = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; Price = '3.621'; Volume = '  -0.34%  ' }
    @{ Row = 34; Coin = 'VeChain'; Link = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; Price = '0.02437'; Volume = '  +0.76%  ' }
    @{ Row = 35; Coin = 'Hedera'; Link = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; Price = '0.06547'; Volume = '  +1.28%  ' }
    @{ Row = 36; Coin = 'FraxShare'; Link = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; Price = '8.941'; Volume = '  +1.08%  ' }
    @{ Row = 37; Coin = 'Algorand'; Link = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; Price = '0.2154'; Volume = '  -1.12%  ' }
    @{ Row = 38; Coin = 'InternetComputer(DFINITY)'; Link = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; Price = '5.012'; Volume = '  +0.16%  ' }
    @{ Row = 39; Coin = 'ARBITRUM'; Link = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; Price = '1.181'; Volume = '  +0.52%  ' }
    @{ Row = 40; Coin = 'TrustWalletToken'; Link = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; Price = '1.234'; Volume = '  -2.69%  ' }
    @{ Row = 41; Coin = 'TheSandbox'; Link = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; Price = '0.6337'; Volume = '  -0.65%  ' }
    @{ Row = 42; Coin = 'Aptos'; Link = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; Price = '11.07'; Volume = '  -1.13%  ' }
    @{ Row = 43; Coin = 'Frax'; Link = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'; Price = '1.006'; Volume = '  -1.07%  ' }
    @{ Row = 44; Coin = 'Decentraland'; Link = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; Price = '0.5964'; Volume = '  -0.62%  ' }
    @{ Row = 45; Coin = 'EnergySwap'; Link = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; Price = '12.95'; Volume = '  +0.11%  ' }
    @{ Row = 46; Coin = 'PancakeSwap'; Link = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; Price = '3.672'; Volume = '  -0.98%  ' }
    @{ Row = 47; Coin = 'NEARProtocol'; Link = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; Price = '1.992'; Volume = '  +0.43%  ' }
    @{ Row = 48; Coin = 'EOS'; Link = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'; Price = '1.207'; Volume = '  +0.54%  ' }
    @{ Row = 49; Coin = 'Quant'; Link = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; Price = '121.71'; Volume = '  -0.10%  ' }
    @{ Row = 50; Coin = 'WEMIXTOKEN'; Link = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; Price = '1.144'; Volume = '  -5.95%  ' }
    @{ Row = 51; Coin = 'Cronos'; Link = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; Price = '0.06786'; Volume = '  -0.66%  ' }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Range("B$r").Value = $item.Coin
    $ws.Range("C$r").Value = $item.Link
    # Price column holds strings such as "28.448.52" that Excel would
    # otherwise auto-coerce into a number (and mangle); force text first.
    $ws.Range("D$r").NumberFormat = "@"
    $ws.Range("D$r").Value = $item.Price
    $ws.Range("D$r").ClearFormats()
    $ws.Range("E$r").Value = $item.Volume
}
